$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.07271233333333334
$ws.Range("H2").Value = 0.218137
$ws.Range("I2").Value = 0.004171225362010892
$ws.Range("J2").Value = 0.004171225362010893
$ws.Range("M2").Value = 0.2901893333333334
$ws.Range("N2").Value = 0.870568
$ws.Range("O2").Value = 0.03429389578125064
$ws.Range("P2").Value = 0.03429389578125064
$ws.Range("Q2").Value = 0.02110034353511111
$ws.Range("R2").Value = 0.189903091816
$ws.Range("S2").Value = 0.000143047567844911
$ws.Range("T2").Value = 0.000143047567844911
$ws.Range("G3").Value = 0.07271233333333334
$ws.Range("H3").Value = 0.218137
$ws.Range("I3").Value = 0.004171225362010892
$ws.Range("J3").Value = 0.004171225362010893
$ws.Range("O3").Value = 0.8402845891331153
$ws.Range("P3").Value = 0.8402845891331153
$ws.Range("Q3").Value = 0.5170101877915556
$ws.Range("R3").Value = 4.653091690124
$ws.Range("S3").Value = 0.003505016389498953
$ws.Range("T3").Value = 0.003505016389498954
$ws.Range("G4").Value = 0.07271233333333334
$ws.Range("H4").Value = 0.218137
$ws.Range("I4").Value = 0.004171225362010892
$ws.Range("J4").Value = 0.004171225362010893
$ws.Range("O4").Value = 0.1254215150856341
$ws.Range("P4").Value = 0.1254215150856341
$ws.Range("Q4").Value = 0.07716933275477779
$ws.Range("R4").Value = 0.6945239947930001
$ws.Range("S4").Value = 0.0005231614046670288
$ws.Range("T4").Value = 0.0005231614046670288
$ws.Range("I5").Value = 0.5387060579248023
$ws.Range("J5").Value = 0.5387060579248023
$ws.Range("M5").Value = 0.2901893333333334
$ws.Range("N5").Value = 0.870568
$ws.Range("O5").Value = 0.03429389578125064
$ws.Range("P5").Value = 0.03429389578125064
$ws.Range("Q5").Value = 2.725070428987556
$ws.Range("R5").Value = 24.525633860888
$ws.Range("S5").Value = 0.01847432940720154
$ws.Range("T5").Value = 0.01847432940720154
$ws.Range("I6").Value = 0.5387060579248023
$ws.Range("J6").Value = 0.5387060579248023
$ws.Range("O6").Value = 0.8402845891331153
$ws.Range("P6").Value = 0.8402845891331153
$ws.Range("S6").Value = 0.4526663985468627
$ws.Range("T6").Value = 0.4526663985468627
$ws.Range("I7").Value = 0.5387060579248023
$ws.Range("J7").Value = 0.5387060579248023
$ws.Range("O7").Value = 0.1254215150856341
$ws.Range("P7").Value = 0.1254215150856341
$ws.Range("Q7").Value = 9.96627691699989
$ws.Range("S7").Value = 0.06756532997073808
$ws.Range("T7").Value = 0.06756532997073807
$ws.Range("I8").Value = 0.4571227167131868
$ws.Range("J8").Value = 0.4571227167131868
$ws.Range("M8").Value = 0.2901893333333334
$ws.Range("N8").Value = 0.870568
$ws.Range("O8").Value = 0.03429389578125064
$ws.Range("P8").Value = 0.03429389578125064
$ws.Range("Q8").Value = 2.312377184938667
$ws.Range("R8").Value = 20.811394664448
$ws.Range("S8").Value = 0.01567651880620419
$ws.Range("T8").Value = 0.01567651880620419
$ws.Range("I9").Value = 0.4571227167131868
$ws.Range("J9").Value = 0.4571227167131868
$ws.Range("O9").Value = 0.8402845891331153
$ws.Range("P9").Value = 0.8402845891331153
$ws.Range("R9").Value = 509.9302315038721
$ws.Range("S9").Value = 0.3841131741967536
$ws.Range("T9").Value = 0.3841131741967536
$ws.Range("I10").Value = 0.4571227167131868
$ws.Range("J10").Value = 0.4571227167131868
$ws.Range("O10").Value = 0.1254215150856341
$ws.Range("P10").Value = 0.1254215150856341
$ws.Range("Q10").Value = 8.456952567722668
$ws.Range("R10").Value = 76.11257310950401
$ws.Range("S10").Value = 0.05733302371022901
$ws.Range("T10").Value = 0.057333023710229
